$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Notes" header in column P
$ws.Range("P1").Value = "Notes"
$ws.Range("P1").Style = $ws.Range("O1").Style

# Add notes to existing row 42 (iter_cPCA run)
$ws.Range("P42").Value = "iter_cPCA"
$ws.Range("P42").Style = $ws.Range("O42").Style

# Insert a new row 43 with data for the "normalize mappedX" experiment
$ws.Range("A43").Value = "ukb51139_subset.csv"
$ws.Range("B43").Value = "28012 x 1081"
$ws.Range("C43").Value = "all"
$ws.Range("D43").Value = "no events"
$ws.Range("E43").Value = "> 140/80"
$ws.Range("F43").Value = "zscore"
$ws.Range("G43").Value = "median"
$ws.Range("H43").Value = "none"
$ws.Range("I43").Value = 25
$ws.Range("K43").Value = "N/A"
$ws.Range("L43").Value = "-256.1 & -67.7"
$ws.Range("M43").Value = "52.6 & 50.5"
$ws.Range("N43").Value = "N/A"
$ws.Range("O43").Value = "N/A"
$ws.Range("P43").Value = "normalize mappedX in each sub batch"

# Copy styles from row 42 to row 43 for consistency
$ws.Range("A43:O43").Style = $ws.Range("A42:O42").Style
$ws.Range("P43").Style = $ws.Range("O2").Style
